$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-51 down to 43-52.
$ws.Rows(42).Insert()

# Populate the newly inserted row 42 with a new weekly price observation
# (same market/product attributes as the row that followed it, new date).
$ws.Cells.Item(42, 1).Value = 7
$ws.Cells.Item(42, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(42, 3).Value = "Ñuble"
$ws.Cells.Item(42, 4).Value = 44460
$ws.Cells.Item(42, 5).Value = 16
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100108
$ws.Cells.Item(42, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(42, 9).Value = 100108002
$ws.Cells.Item(42, 10).Value = "Mango"
$ws.Cells.Item(42, 11).Value = "Sin especificar"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 60
$ws.Cells.Item(42, 14).Value = 8500
$ws.Cells.Item(42, 15).Value = 9000
$ws.Cells.Item(42, 16).Value = 8750
$ws.Cells.Item(42, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(42, 18).Value = "Brasil"
$ws.Cells.Item(42, 19).Value = 2188
$ws.Cells.Item(42, 20).Value = 4
